# Update counts on the "Inscricoes" sheet to reflect newer registration totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Each entry: Row, Inscritos (E), Pagos (F), Inscrições homologadas (H)
# (G - Isenções deferidas - is unchanged for all these rows)
$updates = @(
    @{ Row = 15; E = 109; F = 46; H = 46 },
    @{ Row = 17; E = 66 },
    @{ Row = 18; E = 59 },
    @{ Row = 24; E = 14 },
    @{ Row = 26; E = 16; F = 8; H = 8 },
    @{ Row = 30; E = 2 },
    @{ Row = 34; E = 10 },
    @{ Row = 36; E = 54 },
    @{ Row = 40; E = 9 },
    @{ Row = 49; E = 40 },
    @{ Row = 59; E = 8 },
    @{ Row = 65; E = 18 },
    @{ Row = 67; E = 23 },
    @{ Row = 68; E = 10 },
    @{ Row = 70; E = 19 },
    @{ Row = 72; E = 23; F = 11; H = 11 },
    @{ Row = 74; E = 10 },
    @{ Row = 76; E = 28 },
    @{ Row = 84; E = 2 },
    @{ Row = 87; E = 4; F = 2; H = 2 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Range("F$r").Value = $u.F }
    if ($u.ContainsKey("H")) { $ws.Range("H$r").Value = $u.H }
}

$wb.Save()
